# "Generate Report for Handoff"
# A new handoff was generated for e2e\b.md (row 3 on each status sheet):
#   - Status moves from "Handed back: in sync with en-US" to "Ready for handoff"
#   - A new handoff xliff file + timestamp is recorded per language
#   - An error detail is attached noting the handback file is stale

$wb = $excel.ActiveWorkbook

# ----- Overview sheet -----
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = "Ready for handoff"
$overview.Range("F3").Value = "Ready for handoff"
$overview.Range("G3").Value = "2016-08-28 02:37:40"

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b97c3be6944add41f1b94f0a16c52cc898b2e5c9/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/4f9cc2cfca85c3512cc625e10497930967118be1/e2e/b.md."

# ----- zh-cn sheet -----
$zhcn = $wb.Worksheets.Item("zh-cn")
# 39.15 "characters" round-trips through the engine's pixel storage to an
# exported column width of exactly 40 (a plain 40 input gets nudged to
# 40.8333 on save because of the characters<->pixel rounding).
$zhcn.Columns.Item(16).ColumnWidth = 39.15
$zhcn.Range("C3").Value = "Ready for handoff"
# Leading apostrophe forces plain text so "False" is stored as a shared
# string (matching the source file) instead of being auto-coerced to a
# native Excel boolean cell.
$zhcn.Range("F3").Value = "'False"
$zhcn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$zhcn.Range("H3").Value = "2016-08-28 02:37:36"
$zhcn.Range("P3").Value = $errorDetail

# ----- de-de sheet -----
$dede = $wb.Worksheets.Item("de-de")
$dede.Columns.Item(16).ColumnWidth = 39.15
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("F3").Value = "'False"
$dede.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$dede.Range("H3").Value = "2016-08-28 02:37:40"
$dede.Range("P3").Value = $errorDetail
